$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 'Nice to meet you, we are——'
$ws.Range("J6").Value = ''
$ws.Range("B7").Value = 'Sir......'
$ws.Range("G9").Value = 'disappear'
$ws.Range("J9").Value = 'disappear'
$ws.Range("B13").Value = 'Sir, this young lady’s skills seem quite impressive'
$ws.Range("B14").Value = 'Yao, you’re being too hasty.'
$ws.Range("B15").Value = 'Look at the way she swings the weapon——when it cuts through the air like a feather, it means she’s not applying power correctly.'
$ws.Range("B18").Value = 'I see. You’re amazing——you can spot these details just from a few moves.'
$ws.Range("J18").Value = ''
$ws.Range("B21").Value = 'May I ask your name please?'
$ws.Range("B23").Value = 'I’m Chen, the best martial artist in this entire manor.'
$ws.Range("B24").Value = 'My apologies, you must be the top disciple of the manor.'
$ws.Range("B25").Value = 'Top disciple? That’s only because there are just two guards left in Qingliu Manor now.'
$ws.Range("B29").Value = 'Enough chatting——who are you anyway?'
$ws.Range("B30").Value = 'I’m Judge Dee, and this is my student, Yao.'
$ws.Range("B40").Value = 'Hello!'
$ws.Range("B41").Value = 'Hello——may I ask if you’re a physician?'
$ws.Range("B43").Value = 'I practice medicine in JiuJiang county at the foot of the mountain. I came up a few days ago for a consultation.'
$ws.Range("B44").Value = 'Are you familiar with the Lord?'

$ws.Range("I10").Select()
